# Insert a new weekly price record as row 118 (pushing existing rows 118-140
# down to 119-141), duplicating the pattern of the prior entry (row 117) but
# with an updated date for the new week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("118:118").Insert()

$ws.Range("A118").Value2 = 4
$ws.Range("B118").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C118").Value2 = "Los Lagos"
$ws.Range("D118").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D118").Value2 = 44522
$ws.Range("E118").Value2 = 10
$ws.Range("F118").Value2 = 100112039
$ws.Range("G118").Value2 = "Ciboulette"
$ws.Range("H118").Value2 = "Sin especificar"
$ws.Range("I118").Value2 = "Primera"
$ws.Range("J118").Value2 = 80
$ws.Range("K118").Value2 = 2500
$ws.Range("L118").Value2 = 2500
$ws.Range("M118").Value2 = 2500
$ws.Range("N118").Value2 = "`$/docena de atados"
$ws.Range("O118").Value2 = "Región Metropolitana"
$ws.Range("P118").Value2 = 833
$ws.Range("Q118").Value2 = 3
$ws.Range("R118").Value2 = "Hortaliza"
